$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.115108251571655
$ws.Range("B1").Value = 4.159480571746826
$ws.Range("C1").Value = 4.434149265289307
$ws.Range("D1").Value = 8.109818458557129
$ws.Range("E1").Value = 3.093109607696533
